# Scheduled market-data refresh: update cached Leve profit figures
# (currentAveragePrice* / LevePrice* / LeveProfit* columns, H:N) across
# several crafting-job sheets to the latest pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 597.45715
$ws.Range("J17").Value = 603.2646999999999
$ws.Range("L17").Value = 1809.7941
$ws.Range("N17").Value = -2145.7941
$ws.Range("H53").Value = 391.52173
$ws.Range("I53").Value = 267.53845
$ws.Range("J53").Value = 552.7
$ws.Range("K53").Value = 267.53845
$ws.Range("L53").Value = 552.7
$ws.Range("M53").Value = 369.46155
$ws.Range("N53").Value = -1826.7
$ws.Range("H55").Value = 83687.414
$ws.Range("J55").Value = 391.66666
$ws.Range("L55").Value = 391.66666
$ws.Range("N55").Value = -819.66666
$ws.Range("H100").Value = 8335545
$ws.Range("I100").Value = 12822140
$ws.Range("J100").Value = 3297.1428
$ws.Range("K100").Value = 12822140
$ws.Range("L100").Value = 3297.1428
$ws.Range("M100").Value = -12821599
$ws.Range("N100").Value = -4379.1428
$ws.Range("H103").Value = 978
$ws.Range("I103").Value = 395
$ws.Range("J103").Value = 1561
$ws.Range("K103").Value = 1185
$ws.Range("L103").Value = 4683
$ws.Range("M103").Value = -599
$ws.Range("N103").Value = -5855
$ws.Range("H113").Value = 267426.25
$ws.Range("I113").Value = 355668.34
$ws.Range("J113").Value = 2700
$ws.Range("K113").Value = 355668.34
$ws.Range("L113").Value = 2700
$ws.Range("M113").Value = -352414.34
$ws.Range("N113").Value = -9208
$ws.Range("H132").Value = 287154.12
$ws.Range("I132").Value = 320520.4
$ws.Range("J132").Value = 75834.336
$ws.Range("K132").Value = 961561.2000000001
$ws.Range("L132").Value = 227503.008
$ws.Range("M132").Value = -959031.2000000001
$ws.Range("N132").Value = -232563.008
$ws.Range("H133").Value = 16666.666
$ws.Range("J133").Value = 16666.666
$ws.Range("L133").Value = 16666.666
$ws.Range("N133").Value = -26786.666
$ws.Range("H135").Value = 5966.4546
$ws.Range("I135").Value = 6107.7144
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 54969.4296
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -52434.4296
$ws.Range("N135").Value = -32070
$ws.Range("H137").Value = 333336670
$ws.Range("I137").Value = 333336670
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 1000010010
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -1000007460
$ws.Range("H138").Value = 6264719.5
$ws.Range("I138").Value = 1346770.9
$ws.Range("J138").Value = 8199650.5
$ws.Range("K138").Value = 4040312.7
$ws.Range("L138").Value = 24598951.5
$ws.Range("M138").Value = -4035172.7
$ws.Range("N138").Value = -24609231.5
# Row 137 no longer has a separate HQ-profit column; the value merges into M137.
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3569.8823
$ws.Range("I61").Value = 2621.1
$ws.Range("J61").Value = 4925.2856
$ws.Range("K61").Value = 2621.1
$ws.Range("L61").Value = 4925.2856
$ws.Range("M61").Value = -2409.1
$ws.Range("N61").Value = -5349.2856
$ws.Range("H74").Value = 7749.3
$ws.Range("I74").Value = 2122
$ws.Range("J74").Value = 18200
$ws.Range("K74").Value = 2122
$ws.Range("L74").Value = 18200
$ws.Range("M74").Value = -1248
$ws.Range("N74").Value = -19948
$ws.Range("H77").Value = 7749.3
$ws.Range("I77").Value = 2122
$ws.Range("J77").Value = 18200
$ws.Range("K77").Value = 10610
$ws.Range("L77").Value = 91000
$ws.Range("M77").Value = -6242
$ws.Range("N77").Value = -99736
$ws.Range("H136").Value = 3569.8823
$ws.Range("I136").Value = 2621.1
$ws.Range("J136").Value = 4925.2856
$ws.Range("K136").Value = 7863.299999999999
$ws.Range("L136").Value = 14775.8568
$ws.Range("M136").Value = -5313.299999999999
$ws.Range("N136").Value = -19875.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1209.421
$ws.Range("I31").Value = 1209.421
$ws.Range("K31").Value = 1209.421
$ws.Range("M31").Value = -914.421
$ws.Range("H34").Value = 1209.421
$ws.Range("I34").Value = 1209.421
$ws.Range("K34").Value = 1209.421
$ws.Range("M34").Value = -1007.421
$ws.Range("H132").Value = 3287.1333
$ws.Range("I132").Value = 2246.7
$ws.Range("K132").Value = 6740.099999999999
$ws.Range("M132").Value = -4210.099999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 30.956522
$ws.Range("I12").Value = 12.6875
$ws.Range("K12").Value = 38.0625
$ws.Range("M12").Value = 134.9375
$ws.Range("H97").Value = 1559.2
$ws.Range("I97").Value = 897.5
$ws.Range("J97").Value = 2000.3334
$ws.Range("K97").Value = 2692.5
$ws.Range("L97").Value = 6001.0002
$ws.Range("M97").Value = -2196.5
$ws.Range("N97").Value = -6993.0002
$ws.Range("H122").Value = 697.4286
$ws.Range("I122").Value = 293.5
$ws.Range("J122").Value = 946
$ws.Range("K122").Value = 2641.5
$ws.Range("L122").Value = 8514
$ws.Range("M122").Value = -191.5
$ws.Range("N122").Value = -13414
$ws.Range("H131").Value = 1379.6438
$ws.Range("I131").Value = 500.81818
$ws.Range("J131").Value = 1535.5646
$ws.Range("K131").Value = 1502.45454
$ws.Range("L131").Value = 4606.6938
$ws.Range("M131").Value = 3537.54546
$ws.Range("N131").Value = -14686.6938

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8075.143
$ws.Range("I22").Value = 495
$ws.Range("J22").Value = 18182
$ws.Range("K22").Value = 495
$ws.Range("L22").Value = 18182
$ws.Range("M22").Value = -200
$ws.Range("N22").Value = -18772
$ws.Range("H27").Value = 8075.143
$ws.Range("I27").Value = 495
$ws.Range("J27").Value = 18182
$ws.Range("K27").Value = 495
$ws.Range("L27").Value = 18182
$ws.Range("M27").Value = -388
$ws.Range("N27").Value = -18396
$ws.Range("H122").Value = 3214.75
$ws.Range("I122").Value = 1433.3334
$ws.Range("J122").Value = 3529.1177
$ws.Range("K122").Value = 4300.0002
$ws.Range("L122").Value = 10587.3531
# Row 122 previously had no LeveProfitNQ figure; add the new cell between L and N.
$ws.Range("M122").Value = -1850.0002
$ws.Range("N122").Value = -15487.3531
$ws.Range("H132").Value = 4625.75
$ws.Range("I132").Value = 2701
$ws.Range("J132").Value = 6550.5
$ws.Range("K132").Value = 8103
$ws.Range("L132").Value = 19651.5
$ws.Range("M132").Value = -5573
$ws.Range("N132").Value = -24711.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 46102.5
$ws.Range("J128").Value = 46102.5
$ws.Range("L128").Value = 46102.5
$ws.Range("N128").Value = -56062.5
